$d = $word.ActiveDocument

# The first paragraph in the body currently carries the spacing
# definition (w:line="220" w:lineRule="atLeast"). We insert a brand new
# paragraph right before it, containing the run "来点不同", formatted
# with the same spacing and an eastAsia font hint (as Word does when
# typing CJK text), then inserting that OOXML fragment into the new
# (still empty) leading paragraph.

$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item(1)
$newRange = $newPara.Range

$owx = @'
<?xml version="1.0" encoding="utf-8"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:spacing w:line="220" w:lineRule="atLeast"/>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>来点不同</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

[void]$newRange.InsertXML($owx)
